$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.829.67'
$ws.Range('E2').Value = '  -0.02%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.144.87'
$ws.Range('E3').Value = '  +0.88%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '533.15'
$ws.Range('E5').Value = '  +1.54%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.59'
$ws.Range('E6').Value = '  -0.84%  '

$ws.Range('E7').Value = '  -0.06%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.142.68'
$ws.Range('E8').Value = '  +0.89%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.449'
$ws.Range('E9').Value = '  +2.96%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.23'
$ws.Range('E10').Value = '  -0.90%  '

$ws.Range('E11').Value = '  +0.20%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.400'
$ws.Range('E12').Value = '  +4.32%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.684.98'
$ws.Range('E13').Value = '  +0.93%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.135'
$ws.Range('E14').Value = '  +3.04%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '25.66'
$ws.Range('E15').Value = '  -2.58%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000166'
$ws.Range('E16').Value = '  +0.28%  '

$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '57.912.47'
$ws.Range('E17').Value = '  -0.06%  '

$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.141.00'
$ws.Range('E18').Value = '  +0.87%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.14'
$ws.Range('E19').Value = '  +0.26%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.83'
$ws.Range('E20').Value = '  -0.23%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.99'
$ws.Range('E21').Value = '  -1.06%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '356.54'
$ws.Range('E22').Value = '  +5.50%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  -0.22%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '68.55'
$ws.Range('E24').Value = '  +2.96%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.511'
$ws.Range('E25').Value = '  +0.55%  '

$ws.Range('E26').Value = '  +0.65%  '

$ws.Range('E27').Value = '  +0.05%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0942'
$ws.Range('E28').Value = '  +1.12%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.46'
$ws.Range('E29').Value = '  +3.24%  '

$ws.Range('E30').Value = '  -0.02%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.41'
$ws.Range('E31').Value = '  -2.90%  '

$ws.Range('E32').Value = '  +1.97%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.31'
$ws.Range('E33').Value = '  +1.79%  '

$ws.Range('E34').Value = '  +0.21%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.88'
$ws.Range('E35').Value = '  +5.61%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.41'
$ws.Range('E36').Value = '  +1.92%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.18'
$ws.Range('E37').Value = '  +1.69%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '26.13'
$ws.Range('E38').Value = '  -3.99%  '

$ws.Range('E39').Value = '  -0.78%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0673'
$ws.Range('E40').Value = '  +0.80%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.63'
$ws.Range('E41').Value = '  +10.83%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.11'
$ws.Range('E42').Value = '  +5.10%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.702'
$ws.Range('E43').Value = '  +2.60%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.183.16'
$ws.Range('E44').Value = '  +0.76%  '

$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0271'
$ws.Range('E46').Value = '  +4.84%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.999'
$ws.Range('E47').Value = '  -0.02%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.328.26'
$ws.Range('E48').Value = '  +1.85%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.00'
$ws.Range('E49').Value = '  +2.25%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.06'
$ws.Range('E50').Value = '  +0.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '20.41'
$ws.Range('E51').Value = '  -1.60%  '
